# Update crypto price/volume figures per latest scrape (GitHub Actions run).
#
# Price (col D) and Volume(1h) (col E) are refreshed with the latest scraped
# text. Some new Price strings read as plain numbers (e.g. "32.10", "0.999");
# left alone, Excel would silently coerce them to numeric values and drop the
# trailing/insignificant zero ("32.1"), which does not match the scraped
# display text. For those cells we briefly force Text format before writing
# the value, then restore the default "Normal" style so no stray formatting
# is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.680.55'
$ws.Range('E2').Value = '  +3.86%  '
$ws.Range('D3').Value = '2.257.79'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.529'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.17%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.479'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.10'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.87'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('E13').Value = '  +1.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.57'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('D15').Value = '2.603.49'
$ws.Range('E15').Value = '  +1.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.78%  '
$ws.Range('D17').Value = '2.260.50'
$ws.Range('E17').Value = '  +3.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.759'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.15%  '
$ws.Range('D19').Value = '41.598.56'
$ws.Range('E19').Value = '  +3.86%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.28'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.86%  '
$ws.Range('D21').Value = '0.0₃0902'
$ws.Range('E21').Value = '  +1.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.57'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.42%  '
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('E27').Value = '  +5.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.92'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.48'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '160.31'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '34.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.74%  '
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.14'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0740'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.78%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('E37').Value = '  +2.37%  '
$ws.Range('E38').Value = '  +2.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.59'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.104'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('E41').Value = '  +2.88%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.91'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.70%  '
$ws.Range('D43').Value = '2.054.20'
$ws.Range('E43').Value = '  -0.98%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.44'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('E45').Value = '  +2.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.20'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('E48').Value = '  +6.28%  '
$ws.Range('E49').Value = '  +3.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.95%  '
$ws.Range('E51').Value = '  +2.27%  '
